$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.471.36"
$ws.Range("E2").Value = "  -0.73%  "

# Row 3
$ws.Range("D3").Value = "1.830.80"
$ws.Range("E3").Value = "  +0.96%  "

# Row 4
$ws.Range("D4").Value = "'0.9959"
$ws.Range("E4").Value = "  -0.64%  "

# Row 5
$ws.Range("D5").Value = "'327.62"
$ws.Range("E5").Value = "  -0.32%  "

# Row 6
$ws.Range("D6").Value = "'0.9922"
$ws.Range("E6").Value = "  -0.67%  "

# Row 7
$ws.Range("D7").Value = "'0.4470"
$ws.Range("E7").Value = "  +1.49%  "

# Row 8
$ws.Range("D8").Value = "'0.3793"
$ws.Range("E8").Value = "  -0.19%  "

# Row 9
$ws.Range("D9").Value = "'45.35"
$ws.Range("E9").Value = "  +1.45%  "

# Row 10
$ws.Range("D10").Value = "'0.07785"
$ws.Range("E10").Value = "  +1.01%  "

# Row 11
$ws.Range("D11").Value = "'1.142"
$ws.Range("E11").Value = "  -0.79%  "

# Row 12
$ws.Range("D12").Value = "'22.33"
$ws.Range("E12").Value = "  -2.09%  "

# Row 13
$ws.Range("D13").Value = "'0.9934"
$ws.Range("E13").Value = "  -0.76%  "

# Row 14
$ws.Range("D14").Value = "'6.337"
$ws.Range("E14").Value = "  -0.08%  "

# Row 15
$ws.Range("D15").Value = "'7.555"
$ws.Range("E15").Value = "  -0.55%  "

# Row 16
$ws.Range("D16").Value = "1.825.36"
$ws.Range("E16").Value = "  +0.68%  "

# Row 17
$ws.Range("D17").Value = "'92.51"
$ws.Range("E17").Value = "  +13.38%  "

# Row 18
$ws.Range("D18").Value = "'0.00001086"
$ws.Range("E18").Value = "  -1.02%  "

# Row 19
$ws.Range("D19").Value = "'0.06378"
$ws.Range("E19").Value = "  -5.47%  "

# Row 20
$ws.Range("D20").Value = "'0.9945"
$ws.Range("E20").Value = "  -0.43%  "

# Row 21
$ws.Range("D21").Value = "'17.60"
$ws.Range("E21").Value = "  -1.05%  "

# Row 22
$ws.Range("D22").Value = "'6.369"
$ws.Range("E22").Value = "  +0.51%  "

# Row 23
$ws.Range("D23").Value = "'0.5378"
$ws.Range("E23").Value = "  -1.21%  "

# Row 24
$ws.Range("D24").Value = "28.530.15"
$ws.Range("E24").Value = "  -0.52%  "

# Row 25
$ws.Range("D25").Value = "'11.84"
$ws.Range("E25").Value = "  -0.43%  "

# Row 26
$ws.Range("E26").Value = "  -11.09%  "

# Row 27
$ws.Range("D27").Value = "'21.02"
$ws.Range("E27").Value = "  +1.08%  "

# Row 28
$ws.Range("D28").Value = "'154.37"
$ws.Range("E28").Value = "  +1.00%  "

# Row 29
$ws.Range("D29").Value = "'2.382"
$ws.Range("E29").Value = "  -0.37%  "

# Row 30
$ws.Range("D30").Value = "2.035.64"
$ws.Range("E30").Value = "  +0.79%  "

# Row 31
$ws.Range("D31").Value = "'130.12"
$ws.Range("E31").Value = "  -2.37%  "

# Row 32
$ws.Range("D32").Value = "'1.228"
$ws.Range("E32").Value = "  -4.50%  "

# Row 33
$ws.Range("D33").Value = "'5.876"
$ws.Range("E33").Value = "  -0.02%  "

# Row 34
$ws.Range("D34").Value = "'0.09271"
$ws.Range("E34").Value = "  -0.31%  "

# Row 35
$ws.Range("D35").Value = "'3.667"
$ws.Range("E35").Value = "  -7.64%  "

# Row 36
$ws.Range("D36").Value = "'12.93"
$ws.Range("E36").Value = "  +4.96%  "

# Row 37
$ws.Range("D37").Value = "'0.02363"
$ws.Range("E37").Value = "  +0.60%  "

# Row 38
$ws.Range("D38").Value = "'0.2203"
$ws.Range("E38").Value = "  -3.06%  "

# Row 39
$ws.Range("D39").Value = "'0.6669"
$ws.Range("E39").Value = "  -0.42%  "

# Row 40
$ws.Range("D40").Value = "'5.214"
$ws.Range("E40").Value = "  -0.71%  "

# Row 41
$ws.Range("D41").Value = "'0.06269"
$ws.Range("E41").Value = "  -2.05%  "

# Row 42
$ws.Range("D42").Value = "'1.193"
$ws.Range("E42").Value = "  -1.64%  "

# Row 43
$ws.Range("D43").Value = "'8.100"
$ws.Range("E43").Value = "  -0.89%  "

# Row 44
$ws.Range("D44").Value = "'1.412"
$ws.Range("E44").Value = "  -2.50%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'14.04"

# Row 46
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").Value = "'0.9915"
$ws.Range("E46").Value = "  -0.63%  "

# Row 47
$ws.Range("D47").Value = "'0.6130"
$ws.Range("E47").Value = "  -0.07%  "

# Row 48
$ws.Range("D48").Value = "'3.757"
$ws.Range("E48").Value = "  -1.48%  "

# Row 49
$ws.Range("D49").Value = "'127.71"
$ws.Range("E49").Value = "  -1.10%  "

# Row 50
$ws.Range("D50").Value = "'2.046"
$ws.Range("E50").Value = "  -0.61%  "

# Row 51
$ws.Range("D51").Value = "'79.67"
$ws.Range("E51").Value = "  +1.07%  "
